# Generate Report for Archive
# Update localization status from "Ready for handoff" to "In Translation"
# on the Overview sheet (zh-cn / de-de status columns) and on each
# per-locale sheet's Status column, then resize the now-narrower status
# columns to fit the shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: column E = zh-cn status, column F = de-de status
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# Per-locale sheets: column C = Status
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# Re-fit the status columns now that the text is shorter.
$overview.Columns.Item(5).ColumnWidth = 13.4101845877511
$overview.Columns.Item(6).ColumnWidth = 13.4101845877511
$zhcn.Columns.Item(3).ColumnWidth = 13.4101845877511
$dede.Columns.Item(3).ColumnWidth = 13.4101845877511
